# mz_lf_tas1_2_participant_202011.xlsx : bump survey from v2 -> v3
#  - rename the recorder-id variable p_recorderID -> p_recorder_id
#  - include the recorder id in the generated participant code formula
#  - bump the displayed form title / form_id to V3 / _v3
#  - restore the view/selection state captured in the saved file

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# 1. Rename the recorder-id question name.
$survey.Range("B2").Value = "p_recorder_id"

# 2. Update the calculated participant-code formula to weave in the
#    (renamed) recorder id between the cluster id and the sequence id.
$survey.Range("L12").Value = "concat(${p_cluster_id}, '-', ${p_recorder_id}, '-', ${p_id_sequence})"

# 3. Bump the form title and form id shown on the settings sheet.
$settings.Range("A2").Value = "2. TAS FL - Inscrição V3"
$settings.Range("B2").Value = "mz_lf_tas1_2_participant_202011_v3"

# 4. Restore saved sheet view / selection state.
$survey.Activate()
$survey.Application.ActiveWindow.Panes.Item(4).ScrollColumn = 6
$survey.Range("L12").Select()

$settings.Activate()
$settings.Range("A2").Select()

# Re-activate the tab that was active when the workbook was last saved.
$survey.Activate()
